$wb = $excel.ActiveWorkbook

# The "Repayment schedule" sheet becomes the active/selected tab
# (previously "NewLoanInput" was selected).
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate() | Out-Null

# Insert a new blank column before column N ("Late"), shifting the
# existing N/O/P columns (Late / heading-Date / Outstanding) one to
# the right. The new column inherits column M's width (11 chars) but
# without the bestFit flag, matching a normal Excel column insert.
$ws.Columns("N").Insert() | Out-Null
$ws.Columns("N").ColumnWidth = 10.166666666666666

# Update the remembered selection on the sheet to J20.
$ws.Range("J20").Select() | Out-Null
